$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13, shifting the existing rows 13-16 down to 14-17.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record
# (same categorical columns as the old row 13, new date/price/origin data).
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value = 44540
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103003
$ws.Range("J13").Value = "Damasco"
$ws.Range("K13").Value = "Castle Brite"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 600
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 16000
$ws.Range("Q13").Value = "$/caja 18 kilos"
$ws.Range("R13").Value = "Región del Maule"
$ws.Range("S13").Value = 889
$ws.Range("T13").Value = 18
